$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.401418566703796
$ws.Range("B1").Value = 3.092527151107788
$ws.Range("C1").Value = 3.949459791183472
$ws.Range("D1").Value = 2.758872747421265
$ws.Range("E1").Value = 0.7381110191345215
